$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 6, 7, 9)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "lipid/free"
}
